# Remove the hard-coded Oracle Cloud credentials (URL / user / password)
# from the "Input_Value" sheet before re-uploading the test-data workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")

# Clear the stored URL, UserName and Password values (L2:N2). Clearing the
# values also drops the now-unused shared strings for those three secrets.
$ws.Range("L2:N2").ClearContents()

# Reflect the on-screen selection that was left after editing those cells.
$ws.Activate()
$ws.Range("L2:N2").Select()
